$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Headers: both the "default" and "first page" headers contain the BTec
# logo picture, currently named "image1.jpg" -> should become "image2.jpg"
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hdr = $sec.Headers.Item($i)
    $shapes = $hdr.Range.InlineShapes
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shape = $shapes.Item($j)
        if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
            $shape.Name = "image2.jpg"
        }
    }
}

# Footers: both the "default" and "first page" footers contain the Pearson
# Edexcel logo picture, currently named "image2.png" -> should become "image1.png"
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $ftr = $sec.Footers.Item($i)
    $shapes = $ftr.Range.InlineShapes
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shape = $shapes.Item($j)
        if ($shape.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shape.Name = "image1.png"
        }
    }
}
